$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")
$ws.Range("C12").Value = "Yes"
$ws.Range("C14").Value = "Yes"
